# "Generate Report for Handoff"
# A new handoff-xliff generation pass ran for the file
# 38fdfe28-5d2b-49c5-9a5c-d94814134381.md (row 6 in every sheet) and
# refreshed its "Latest Handoff / HO Xliff Generate Date" timestamps,
# which previously (incorrectly) still showed the same datetime as the
# row above it.

$wb = $excel.ActiveWorkbook

# Overview sheet: column G = "Latest HO Xliff Generate Date", row 6
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G6").Value = "2016-08-16 12:41:13"

# zh-cn sheet: column H = "Latest Handoff Datetime", row 6
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H6").Value = "2016-08-16 12:41:03"

# de-de sheet: column H = "Latest Handoff Datetime", row 6
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H6").Value = "2016-08-16 12:41:13"
